$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.903.45'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.552.27'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("E4").Value = '  +0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.73'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.486'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.68'
$ws.Range("E8").Value = '  +1.55%  '
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.773.74'
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.553.77'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.904.73'
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '216.93'
$ws.Range("E18").Value = '  +2.05%  '
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("E23").Value = '  +1.26%  '
$ws.Range("E24").Value = '  +1.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.68'
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.88'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("E30").Value = '  +2.82%  '
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.422.24'
$ws.Range("E33").Value = '  +4.91%  '
$ws.Range("E34").Value = '  +3.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +4.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.958'
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0164'
$ws.Range("E38").Value = '  +1.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.522'
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.807'
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.65'
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("E43").Value = '  -0.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.27'
$ws.Range("E44").Value = '  +4.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.61'
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.688.25'
$ws.Range("E47").Value = '  +1.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.22'
$ws.Range("E48").Value = '  +0.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0520'
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  +5.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0955'
$ws.Range("E51").Value = '  +1.24%  '
